$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [string][char]0x00A0

# --- Insert 4 new rows at the top (rows 1-4), pushing existing data down ---
$ws.Rows("1:4").Insert()

# Carry over the number formats from the (now shifted) first data row, so the
# new rows pick up the same cell styles (date / text / text-left) that the
# rest of the sheet uses instead of creating brand-new style entries.
$ws.Range("A5:G5").Copy()
$ws.Range("A1:G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 1 : newest movement ---
$ws.Range("A1").Value = 41680
$ws.Range("B1").Value = $nbsp + $nbsp + "TRANSFERENCIA INTERNET"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "0000952785"
$ws.Range("E1").Value = "AG. NORTE"
$ws.Range("F1").Value = "100.00" + $nbsp + $nbsp
$ws.Range("G1").Value = "2497.81"

# --- Row 2 ---
$ws.Range("A2").Value = 41677
$ws.Range("B2").Value = "DEPOSITO"
$ws.Range("C2").Value = "C"
$ws.Range("D2").Value = "0002420226"
$ws.Range("E2").Value = "TENA"
$ws.Range("F2").Value = "100.00" + $nbsp + $nbsp
$ws.Range("G2").Value = "2597.81"

# --- Row 3 ---
$ws.Range("A3").Value = 41677
$ws.Range("B3").Value = "PAGO PRESTAMO"
$ws.Range("C3").Value = "D"
$ws.Range("D3").Value = "0000937060"
$ws.Range("E3").Value = "AG. NORTE"
$ws.Range("F3").Value = "281.05" + $nbsp + $nbsp
$ws.Range("G3").Value = "2497.81"

# --- Row 4 ---
$ws.Range("A4").Value = 41676
$ws.Range("B4").Value = $nbsp + $nbsp + "TRANSFERENCIA INTERNET"
$ws.Range("C4").Value = "C"
$ws.Range("D4").Value = "0004190618"
$ws.Range("E4").Value = "AG. NORTE"
$ws.Range("F4").Value = "281.05" + $nbsp + $nbsp
$ws.Range("G4").Value = "2778.86"

# --- H1 formula: now also stamps the real creation timestamp via NOW() ---
$ws.Range("H1").Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B1,""', 'mo_tipo' => '"",C1,""', 'mo_documento' => '"",D1,""', 'mo_oficina' => '"",E1,""', 'mo_monto' => "",TRIM(F1),"", 'mo_saldo' => "",G1,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL),"")"

# --- Sheet view: park the selection on I10 (matches the saved view state) ---
[void]$ws.Range("I10").Select()
